$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "CR0 Subgroup 1"
$ws.Range("B2").Value = 2633.265520811323
$ws.Range("C2").Value = 0.0007202279044656574
$ws.Range("D2").Value = 0.05450764348867843

$ws.Range("A3").Value = "CR1 Subgroup 1"
$ws.Range("B3").Value = 2741.885299681539
$ws.Range("C3").Value = 0.0007906561501348794
$ws.Range("D3").Value = 0.05825222639085226

$ws.Range("A4").Value = "CR2 Subgroup 1"
$ws.Range("B4").Value = 2907.059354896189
$ws.Range("C4").Value = 0.002725059735176096
$ws.Range("D4").Value = 0.064063329474777

$ws.Range("A5").Value = "CR3 Subgroup 1"
$ws.Range("B5").Value = 2962.715762636106
$ws.Range("C5").Value = 0.0009510932762169244
$ws.Range("D5").Value = 0.06605056052768943

$ws.Range("A6").Value = "DM1 Subgroup 1"
$ws.Range("B6").Value = 11991.56429097778
$ws.Range("C6").Value = 0.001855050517804356
$ws.Range("D6").Value = 0.2356483625798341

$ws.Range("A7").Value = "DM2 Subgroup 1"
$ws.Range("B7").Value = 12063.57947437494
$ws.Range("C7").Value = 0.0009664084098867688
$ws.Range("D7").Value = 0.2371782519274735

$ws.Range("A8").Value = "DM3 Subgroup 1"
$ws.Range("B8").Value = 12087.76090545364
$ws.Range("C8").Value = 0.0002177353319765159
$ws.Range("D8").Value = 0.237691382819225

$ws.Range("A9").Value = "NG1.5 Subgroup 1"
$ws.Range("B9").Value = -413.7542166640237
$ws.Range("C9").Value = -0.0001417506002811621
$ws.Range("D9").Value = 0.003063051899264035

$ws.Range("A10").Value = "NG2 Subgroup 1"
$ws.Range("B10").Value = -1991.573974836802
$ws.Range("C10").Value = -0.0007437054569524939
$ws.Range("D10").Value = 0.1092794004752168

$ws.Range("A11").Value = "NG3 Subgroup 1"
$ws.Range("B11").Value = -3569.393733009581
$ws.Range("C11").Value = -0.0014528508583018
$ws.Range("D11").Value = 0.4545533609356666

$ws.Range("A12").Value = "CR0 Subgroup 2"
$ws.Range("B12").Value = 88862.50469132305
$ws.Range("C12").Value = 3272643.68897705
$ws.Range("D12").Value = 0.878690686980199

$ws.Range("A13").Value = "CR1 Subgroup 2"
$ws.Range("B13").Value = 117187.7235918807
$ws.Range("C13").Value = 2007404786.309897
$ws.Range("D13").Value = 0.911050929229521

$ws.Range("A14").Value = "CR2 Subgroup 2"
$ws.Range("B14").Value = 176932.3201075681
$ws.Range("C14").Value = 3678255467297659
$ws.Range("D14").Value = 0.9373314311454858

$ws.Range("A15").Value = "CR3 Subgroup 2"
$ws.Range("B15").Value = 200808.2932331476
$ws.Range("C15").Value = 240594460968619200
$ws.Range("D15").Value = 0.9408074290888392

$ws.Range("A16").Value = "DM1 Subgroup 2"
$ws.Range("B16").Value = 187554.2156777986
$ws.Range("C16").Value = 1094853906542244
$ws.Range("D16").Value = 0.8901548007476798

$ws.Range("A17").Value = "DM2 Subgroup 2"
$ws.Range("B17").Value = 203593.8047787519
$ws.Range("C17").Value = 19394638841616620
$ws.Range("D17").Value = 0.8993983196477597

$ws.Range("A18").Value = "DM3 Subgroup 2"
$ws.Range("B18").Value = 210086.8295971398
$ws.Range("C18").Value = 18172263492392380
$ws.Range("D18").Value = 0.903029870512559

$ws.Range("A19").Value = "NG1.5 Subgroup 2"
$ws.Range("B19").Value = 74848.74696286963
$ws.Range("C19").Value = 305771.3289966384
$ws.Range("D19").Value = 0.9035493940430257

$ws.Range("A20").Value = "NG2 Subgroup 2"
$ws.Range("B20").Value = 53679.25864836411
$ws.Range("C20").Value = 3386.4793407898
$ws.Range("D20").Value = 0.8951279902880382

$ws.Range("A21").Value = "NG3 Subgroup 2"
$ws.Range("B21").Value = 32509.77033385858
$ws.Range("C21").Value = 31.6727313701562
$ws.Range("D21").Value = 0.8749082534798314

$ws.Range("A22").Value = "CR0 Subgroup 3"
$ws.Range("B22").Value = -12038.50675612018
$ws.Range("C22").Value = -0.001702427197863892
$ws.Range("D22").Value = 0.9617670625144055

$ws.Range("A23").Value = "CR1 Subgroup 3"
$ws.Range("B23").Value = 10260.55743790905
$ws.Range("C23").Value = 0.1886857571386683
$ws.Range("D23").Value = 0.4708667833885153

$ws.Range("A24").Value = "CR2 Subgroup 3"
$ws.Range("B24").Value = 165147.7161863035
$ws.Range("C24").Value = 48143306173680.88
$ws.Range("D24").Value = 0.8597812835487186

$ws.Range("A25").Value = "CR3 Subgroup 3"
$ws.Range("B25").Value = 224194.7135616544
$ws.Range("C25").Value = 1646747813823699000
$ws.Range("D25").Value = 0.8649734851539975

$ws.Range("A26").Value = "DM1 Subgroup 3"
$ws.Range("B26").Value = -9839.754331927435
$ws.Range("C26").Value = -0.001822410866417914
$ws.Range("D26").Value = 0.8721347581126049

$ws.Range("A27").Value = "DM2 Subgroup 3"
$ws.Range("B27").Value = -6342.693906997762
$ws.Range("C27").Value = -0.001797204070015411
$ws.Range("D27").Value = 0.5582135757219189

$ws.Range("A28").Value = "DM3 Subgroup 3"
$ws.Range("B28").Value = -3082.706273283913
$ws.Range("C28").Value = -0.0004292508214496291
$ws.Range("D28").Value = 0.1575914668306114

$ws.Range("A29").Value = "NG1.5 Subgroup 3"
$ws.Range("B29").Value = 2094.61856516839
$ws.Range("C29").Value = 0.006949074947127502
$ws.Range("D29").Value = 0.07190325599109482

$ws.Range("A30").Value = "NG2 Subgroup 3"
$ws.Range("B30").Value = -1988.350871201937
$ws.Range("C30").Value = -0.002801831497230764
$ws.Range("D30").Value = 0.1034152682097278

$ws.Range("A31").Value = "NG3 Subgroup 3"
$ws.Range("B31").Value = -6071.320307572264
$ws.Range("C31").Value = -0.003633783738632058
$ws.Range("D31").Value = 0.6771229474814363
